$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.972.31"
$ws.Range("E2").Value = "  -1.51%  "

$ws.Range("D3").Value = "2.645.18"
$ws.Range("E3").Value = "  -0.68%  "

$ws.Range("E4").Value = "  +0.11%  "

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.53"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +0.61%  "

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.57"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -0.83%  "

$ws.Range("E7").Value = "  +0.27%  "

$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.572"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  -1.14%  "

$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.95"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  +9.68%  "

$ws.Range("E10").Value = "  -2.95%  "

$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.336"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  -0.94%  "

$ws.Range("E12").Value = "  +1.91%  "

$ws.Range("D13").Value = "3.111.71"
$ws.Range("E13").Value = "  -0.46%  "

$ws.Range("D14").Value = "59.021.64"
$ws.Range("E14").Value = "  -1.39%  "

$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.19"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  -0.34%  "

$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000136"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  -1.66%  "

$ws.Range("D17").Value = "2.648.47"
$ws.Range("E17").Value = "  -0.91%  "

$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "340.58"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  -2.99%  "

$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.41"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  -2.64%  "

$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.35"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  -0.70%  "

$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.35"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  +1.08%  "

$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  -0.28%  "

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "63.67"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +1.28%  "

$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.413"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  -1.31%  "

$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.167"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  +0.22%  "

$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  +0.62%  "

$ws.Range("D27").Value = "0.0₃0803"
$ws.Range("E27").Value = "  -1.43%  "

$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.12"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  -0.90%  "

$ws.Range("E29").Value = "  +0.55%  "

$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  +0.04%  "

$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.59"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +0.09%  "

$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.79"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  -0.87%  "

$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.00"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  -1.02%  "

$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.18"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +1.84%  "

$ws.Range("E35").Value = "  -0.03%  "

$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.902"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  -4.65%  "

$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.882"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  +0.64%  "

$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.65"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  -0.54%  "

$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.48"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  -3.59%  "

$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.61"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  -2.63%  "

$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.618"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +1.75%  "

$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  +0.38%  "

$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "275.89"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  -2.49%  "

$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.88"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  +0.35%  "

$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0972"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  -1.65%  "

$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0539"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  +0.77%  "

$ws.Range("D47").Value = "2.048.50"
$ws.Range("E47").Value = "  -1.94%  "

$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.52"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  +2.18%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.79"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  +0.92%  "

$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.09"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  +0.13%  "

$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0228"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  -2.07%  "
